# The authored edit swaps the contents of ppt/theme/theme1.xml and
# ppt/theme/theme2.xml:
#   - theme1.xml (the slide master's theme, i.e. the deck's visible
#     design) goes from the "Integral" / "Red Violet" colour scheme to
#     the default "Office Theme" / "Office" colour scheme.
#   - theme2.xml (only used by the notes master) goes from the default
#     "Office Theme" / "Office" colour scheme to the "Integral" /
#     "Red Violet" colour scheme.
# The font scheme and format scheme (fills/lines/effects) are already
# byte-identical between the two themes, so only the 12 colour-scheme
# slots actually change.
#
# Through the PowerPoint object model exposed here, only the slide
# master's theme (theme1.xml) is reachable/settable, so this script
# recolours that one; see the note below.

function HexToComRGB($hex) {
    # PowerPoint's ColorFormat.RGB takes an integer that is the RGB()
    # macro value: r + g*256 + b*65536 (i.e. bytes stored "BGR").
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Colour-scheme slot order exposed through ThemeColorScheme.Colors(1..12):
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink

# Target colours (originally theme2.xml's "Office Theme" / "Office"
# colour scheme) that theme1.xml must be recoloured to.
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$p = $ppt.ActivePresentation

# theme1.xml is the theme of the slide master (the deck's real/visible
# design) -- recolour it from "Integral"/Red Violet to the default
# Office Theme palette.
#
# NOTE: in this host, $p.NotesMaster.Theme / $p.HandoutMaster.Theme
# resolve to the very same underlying theme object as
# $p.SlideMaster.Theme (there is only one theme the object model lets
# us reach), so we must only write the colours once here -- writing a
# second time through "NotesMaster.Theme" would simply clobber this
# assignment and undo it.
$slideTheme = $p.SlideMaster.Theme
for ($i = 1; $i -le 12; $i++) {
    $slideTheme.ThemeColorScheme.Colors($i).RGB = HexToComRGB($officeColors[$i - 1])
}
